$wb = $excel.ActiveWorkbook

# --- 1. "总计" sheet: make room for a new 2022-Q4 row right under the header,
#        by copying the existing data block down one row (keeps formatting) ---
$summary = $wb.Worksheets.Item(1)
$summary.Range("A2:D6").Copy($summary.Range("A3:D7"))

# Fill in the new 2022-Q4 row
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 1
$summary.Cells.Item(2,4).Value = 0.68

# Re-sequence the index column (A) for the rows that moved down
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5

# --- 2. Add the new "2022-Q4" worksheet right after "总计" by duplicating the
#        "2022-Q3" sheet (so formatting / sheet properties match its siblings) ---
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $summary)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

$newSheet.Cells.Item(2,2).Value = "010695"
$newSheet.Cells.Item(2,3).Value = "华夏磐益一年定期开放混合"

# D/E/F/G hold numeric-looking figures that are stored as *text* in this
# workbook. A leading apostrophe keeps them text; re-applying the "Normal"
# style afterwards strips the quote-prefix formatting Excel would otherwise
# stamp on the cell, so the XML stays byte-for-byte like its siblings.
$newSheet.Cells.Item(2,4).Value = "'16.03"
$newSheet.Cells.Item(2,4).Style = "Normal"
$newSheet.Cells.Item(2,5).Value = "'98.69"
$newSheet.Cells.Item(2,5).Style = "Normal"
$newSheet.Cells.Item(2,6).Value = "'4.25"
$newSheet.Cells.Item(2,6).Style = "Normal"
$newSheet.Cells.Item(2,7).Value = "'0.6813"
$newSheet.Cells.Item(2,7).Style = "Normal"

$newSheet.Cells.Item(2,8).Value = 8

# Restore the originally-selected tab (2021-Q3, the last sheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
